$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.783.42"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.15%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.861.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.82%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6413"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.46%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.76"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07552"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.2975"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.26%  "
$ws.Range("E11").Value = "  +5.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07677"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.863.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.81%  "
$ws.Range("E14").Value = "  +1.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6934"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "84.07"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000009884"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +10.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.118"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "29.788.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.113.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "236.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.506"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.001"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1423"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.558"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "17.94"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06175"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.492"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.290"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.164"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.103"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.897"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.176"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.7311"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.92%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.603"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.827"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01786"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.206.65"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9212"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.286"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.08%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.021.35"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.02"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "66.62"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.62%  "
$ws.Range("E48").Value = "  +1.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4065"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.180"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.669"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.39%  "
